$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "22.359.01"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -0.07%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.566.78"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +0.14%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.004"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  +0.32%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "1.003"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.19%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "290.57"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.01%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.3741"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +0.92%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "49.02"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -0.40%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.3382"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -0.09%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.07527"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -1.54%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "1.129"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -3.07%  "

$ws.Range("E12").Value = "  +0.31%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "20.83"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -3.05%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "5.927"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -1.86%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "6.878"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -0.62%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "1.565.28"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -0.32%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.00001119"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -0.72%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "89.62"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -0.59%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.06732"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +0.16%  "

$ws.Range("E20").Value = "  +0.27%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "6.170"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -1.15%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "16.42"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -0.68%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "11.88"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -1.19%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "22.357.16"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -0.06%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.378"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +0.39%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "2.707"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -3.52%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "20.02"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -0.61%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "147.67"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +1.68%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "5.039"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +1.10%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "125.16"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -0.33%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "1.739.99"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +0.37%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "2.018"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +0.59%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.9840"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -1.74%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "5.999"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -3.56%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "9.974"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -1.22%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "1.409"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +8.84%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.08461"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -0.22%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.02473"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -2.37%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.2275"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -2.05%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.06437"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +0.31%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "5.373"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -2.78%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.6246"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -1.38%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "11.00"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -6.01%  "

$ws.Range("E44").Value = "  +0.20%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "13.96"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -1.25%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "3.796"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +0.98%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.5874"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -1.79%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "2.054"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -2.06%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "1.254"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -0.84%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "124.24"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -0.11%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.07319"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +0.66%  "
